$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.439.53"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "2.918.84"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.67%  "

$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.17%  "

$ws.Range("D14").Value = "3.370.90"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("D16").Value = "2.911.87"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.958"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.57%  "

$ws.Range("D18").Value = "51.367.06"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.95%  "

$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.174"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.78%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.24%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.27%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.77%  "

$ws.Range("E42").Value = "  -1.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.37%  "

$ws.Range("E45").Value = "  -3.15%  "

$ws.Range("D46").Value = "2.092.28"
$ws.Range("E46").Value = "  -4.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.238"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0335"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.903"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.13%  "
